# Generate Report for Handback
# - Mark the two handoff rows as "Handed back: in sync with en-US" on the
#   Overview sheet as well as the per-language (zh-cn / de-de) detail sheets.
# - Populate the "Latest Target File" / "Latest Handback File" columns (F/G)
#   for each detail sheet, with hyperlinks matching the source/handback files.
# - Stamp the handback datetime into the "Latest Handback DateTime" column (H).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile   = "b0c057c4-3620-41cb-8935-787f8ae38e7b.md"
$zhXlf    = "b0c057c4-3620-41cb-8935-787f8ae38e7b.65970977ac92851fd0dba10b197b1fb69e45c616.zh-cn.xlf"
$deXlf    = "b0c057c4-3620-41cb-8935-787f8ae38e7b.65970977ac92851fd0dba10b197b1fb69e45c616.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/8e6d5876ecad64f3d8353e7151180e906a313b94/e2e/b0c057c4-3620-41cb-8935-787f8ae38e7b.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac52152e0af90e74eefdb2c06cee3aa22de8cc03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b0c057c4-3620-41cb-8935-787f8ae38e7b.65970977ac92851fd0dba10b197b1fb69e45c616.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/055f83b20e9ca776e93ca20770539e4f53de6b70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b0c057c4-3620-41cb-8935-787f8ae38e7b.65970977ac92851fd0dba10b197b1fb69e45c616.de-de.xlf"

# ---- Overview sheet: status columns for both rows ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---- zh-cn detail sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, $null, $null, $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, $null, $null, $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, $null, $null, $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, $null, $null, $zhXlf)

$wsZh.Range("H2").Value = "2016-03-18 10:47:44"
$wsZh.Range("H3").Value = "2016-03-18 10:47:44"

# ---- de-de detail sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, $null, $null, $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, $null, $null, $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, $null, $null, $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, $null, $null, $deXlf)

$wsDe.Range("H2").Value = "2016-03-18 10:47:49"
$wsDe.Range("H3").Value = "2016-03-18 10:47:49"
